$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 63 - this pushes the existing rows
# 63-67 down to 65-69, matching the target dimension A1:R69.
$ws.Rows("63:64").Insert()

# New row 63: Cilantro "Primera" record dated 44782
$ws.Range("A63").Value = 7
$ws.Range("B63").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C63").Value = "Ñuble"
$ws.Range("D63").Value = 44782
$ws.Range("E63").Value = 16
$ws.Range("F63").Value = 100112040
$ws.Range("G63").Value = "Cilantro"
$ws.Range("H63").Value = "Sin especificar"
$ws.Range("I63").Value = "Primera"
$ws.Range("J63").Value = 200
$ws.Range("K63").Value = 700
$ws.Range("L63").Value = 800
$ws.Range("M63").Value = 750
$ws.Range("N63").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O63").Value = "Provincia de Diguillín"
$ws.Range("P63").Value = 750
$ws.Range("Q63").Value = 1
$ws.Range("R63").Value = "Hortaliza"

# New row 64: Cilantro "Segunda" record dated 44782
$ws.Range("A64").Value = 7
$ws.Range("B64").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C64").Value = "Ñuble"
$ws.Range("D64").Value = 44782
$ws.Range("E64").Value = 16
$ws.Range("F64").Value = 100112040
$ws.Range("G64").Value = "Cilantro"
$ws.Range("H64").Value = "Sin especificar"
$ws.Range("I64").Value = "Segunda"
$ws.Range("J64").Value = 150
$ws.Range("K64").Value = 600
$ws.Range("L64").Value = 600
$ws.Range("M64").Value = 600
$ws.Range("N64").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O64").Value = "Provincia de Diguillín"
$ws.Range("P64").Value = 600
$ws.Range("Q64").Value = 1
$ws.Range("R64").Value = "Hortaliza"
